# Soil and Leaf Analysis SVS October 2016.xlsx - apply commit changes
#
# Summary of the edit (per commit message / xml diff):
#  - Sheet3 (previously blank) is populated with the full list of sample
#    names (a copy of the "Soil" sheet's sample-name column, B4:B117),
#    with a bold "Sample" header in A1 and the data following in A2:A114.
#  - Column A on Sheet3 is sized to fit its (now much longer) contents.
#  - The "Soil" sheet is no longer the tab shown when the workbook opens;
#    instead its selection becomes the whole of column B, and Sheet3
#    becomes the selected/active tab, scrolled down a bit with B11 selected.
#  - (Workbook-level bookkeeping - calcId, fileVersion, mc:AlternateContent,
#    etc. - is maintained automatically by the host application on save.)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets("Soil")
$ws3 = $wb.Worksheets("Sheet3")

# --- Populate Sheet3 with the sample-name list copied from Soil!B4:B117 ---
$names = $ws1.Range("B4:B117").Value2
$ws3.Range("A1:A114").Value2 = $names

# Bold the header cell (A1, "Sample") to match the style used for the
# equivalent header cell on the other sheets.
$ws3.Range("A1").Font.Bold = $true

# Size column A to fit the (now populated) contents.
$ws3.Columns.Item(1).AutoFit()

# --- Update view/selection state ---
# "Soil" is no longer the selected tab; its selection becomes column B.
$ws1.Activate()
$ws1.Columns.Item(2).Select()

# Sheet3 becomes the active/selected tab, scrolled so row 7 is at the top,
# with B11 as the active selected cell.
$ws3.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws3.Range("B11").Select()
